$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the 2025-08-24 buy as a new row right after the current last row (52).
# Column A stores the date as literal text (matching the existing rows'
# "MM/DD/2025" inline-string convention), so force a text format before
# assigning the value to stop Excel from auto-converting it to a date
# serial number. Reset the style back to Normal afterwards so the cell
# doesn't carry a stray number-format override (matching the unstyled
# A42:A52 cells above it).
$row = 53
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "08/24/2025"
$dateCell.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 112.0649999999951
$ws.Cells.Item($row, 3).Value = 0.08923392673895009
$ws.Cells.Item($row, 4).Value = 10
